$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.014401056327781
$ws.Range("D2").Value = 1.021002681843134
$ws.Range("E2").Value = 1.016198192904378
$ws.Range("I2").Value = 1.026041344343634
$ws.Range("J2").Value = 1.019632505049139
$ws.Range("K2").Value = 1.023841298270967
$ws.Range("L2").Value = 1.019051080092922
$ws.Range("N2").Value = 1.010676914312065
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.015190588702211
$ws.Range("D3").Value = 1.021562880961042
$ws.Range("E3").Value = 1.016863031983433
$ws.Range("I3").Value = 1.02611002221871
$ws.Range("J3").Value = 1.02005744383221
$ws.Range("K3").Value = 1.02420853212053
$ws.Range("L3").Value = 1.019521635748569
$ws.Range("N3").Value = 1.010817989540165
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.015702076021302
$ws.Range("D4").Value = 1.021925588425563
$ws.Range("E4").Value = 1.017294161022259
$ws.Range("I4").Value = 1.026153146442554
$ws.Range("J4").Value = 1.020332368821031
$ws.Range("K4").Value = 1.02444567480292
$ws.Range("L4").Value = 1.019826369861695
$ws.Range("N4").Value = 1.010909235335247
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.015917248766061
$ws.Range("D5").Value = 1.022078121572416
$ws.Range("E5").Value = 1.017475629210885
$ws.Range("I5").Value = 1.026170960337363
$ws.Range("J5").Value = 1.02044793660265
$ws.Range("K5").Value = 1.024545252674729
$ws.Range("L5").Value = 1.019954539184813
$ws.Range("N5").Value = 1.010947585166571
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.015953385576144
$ws.Range("D6").Value = 1.022103735486561
$ws.Range("E6").Value = 1.017506111451608
$ws.Range("I6").Value = 1.026173932844794
$ws.Range("J6").Value = 1.020467340274742
$ws.Range("K6").Value = 1.024561965338757
$ws.Range("L6").Value = 1.019976062800238
$ws.Range("N6").Value = 1.010954023677141
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.015704950608619
$ws.Range("D7").Value = 1.021927626383374
$ws.Range("E7").Value = 1.017296584941081
$ws.Range("I7").Value = 1.02615338571341
$ws.Range("J7").Value = 1.020333913087115
$ws.Range("K7").Value = 1.024447005828564
$ws.Range("L7").Value = 1.019828082236123
$ws.Range("N7").Value = 1.010909747807091
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.014667755672212
$ws.Range("D8").Value = 1.021191956748328
$ws.Range("E8").Value = 1.016422683861095
$ws.Range("I8").Value = 1.026064826186389
$ws.Range("J8").Value = 1.01977612230275
$ws.Range("K8").Value = 1.023965505619294
$ws.Range("L8").Value = 1.019210053264768
$ws.Range("N8").Value = 1.010724599191293
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.012844817457269
$ws.Range("D9").Value = 1.019897396021618
$ws.Range("E9").Value = 1.014889997361929
$ws.Range("I9").Value = 1.025898734105853
$ws.Range("J9").Value = 1.018792987791979
$ws.Range("K9").Value = 1.023113413127184
$ws.Range("L9").Value = 1.018123012253631
$ws.Range("N9").Value = 1.010398064228391
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.011632817673133
$ws.Range("D10").Value = 1.019035675026208
$ws.Range("E10").Value = 1.013873192135316
$ws.Range("I10").Value = 1.02578129980991
$ws.Range("J10").Value = 1.01813748913099
$ws.Range("K10").Value = 1.022543000637225
$ws.Range("L10").Value = 1.017399755815962
$ws.Range("N10").Value = 1.010180215573127
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.011108812709223
$ws.Range("D11").Value = 1.018662878829981
$ws.Range("E11").Value = 1.013434111981779
$ws.Range("I11").Value = 1.025728867975672
$ws.Range("J11").Value = 1.017853650566721
$ws.Range("K11").Value = 1.022295466792827
$ws.Range("L11").Value = 1.017086937715159
$ws.Range("N11").Value = 1.010085853238142
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.01091429585641
$ws.Range("D12").Value = 1.018524458081827
$ws.Range("E12").Value = 1.013271200887046
$ws.Range("I12").Value = 1.025709155501414
$ws.Range("J12").Value = 1.01774822120502
$ws.Range("K12").Value = 1.022203441855488
$ws.Range("L12").Value = 1.016970798273665
$ws.Range("N12").Value = 1.010050798492327
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.010956014834051
$ws.Range("D13").Value = 1.018554147399618
$ws.Range("E13").Value = 1.013306137563762
$ws.Range("I13").Value = 1.025713394605244
$ws.Range("J13").Value = 1.017770836096673
$ws.Range("K13").Value = 1.022223185106024
$ws.Range("L13").Value = 1.016995708058732
$ws.Range("N13").Value = 1.010058318044192
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.011092731387381
$ws.Range("D14").Value = 1.018651435847692
$ws.Range("E14").Value = 1.013420641958445
$ws.Range("I14").Value = 1.025727243366065
$ws.Range("J14").Value = 1.017844935715858
$ws.Range("K14").Value = 1.02228786160295
$ws.Range("L14").Value = 1.017077336459192
$ws.Range("N14").Value = 1.010082955688246
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.011176983180523
$ws.Range("D15").Value = 1.018711385471383
$ws.Range("E15").Value = 1.013491216182491
$ws.Range("I15").Value = 1.025735744667241
$ws.Range("J15").Value = 1.017890591100967
$ws.Range("K15").Value = 1.022327700391846
$ws.Range("L15").Value = 1.017127637761912
$ws.Range("N15").Value = 1.010098135189321
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.011667611150914
$ws.Range("D16").Value = 1.019060423510655
$ws.Range("E16").Value = 1.013902357964936
$ws.Range("I16").Value = 1.025784746268637
$ws.Range("J16").Value = 1.01815632661284
$ws.Range("K16").Value = 1.022559417363504
$ws.Range("L16").Value = 1.017420524169866
$ws.Range("N16").Value = 1.010186477449894
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.011975584347763
$ws.Range("D17").Value = 1.019279456922195
$ws.Range("E17").Value = 1.014160579930456
$ws.Range("I17").Value = 1.025815060640002
$ws.Range("J17").Value = 1.018323015591788
$ws.Range("K17").Value = 1.022704623393248
$ws.Range("L17").Value = 1.017604340663426
$ws.Range("N17").Value = 1.010241883886288
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.012155296862359
$ws.Range("D18").Value = 1.019407247558401
$ws.Range("E18").Value = 1.014311312389896
$ws.Range("I18").Value = 1.025832589810573
$ws.Range("J18").Value = 1.018420242037937
$ws.Range("K18").Value = 1.022789267268815
$ws.Range("L18").Value = 1.017711592028518
$ws.Range("N18").Value = 1.01027419836083
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.012216587181795
$ws.Range("D19").Value = 1.019450826253408
$ws.Range("E19").Value = 1.014362727868828
$ws.Range("I19").Value = 1.025838540883989
$ws.Range("J19").Value = 1.018453393629946
$ws.Range("K19").Value = 1.022818119718585
$ws.Range("L19").Value = 1.017748167745866
$ws.Range("N19").Value = 1.01028521620997
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.011942533764466
$ws.Range("D20").Value = 1.019255953369334
$ws.Range("E20").Value = 1.014132863148683
$ws.Range("I20").Value = 1.025811823983013
$ws.Range("J20").Value = 1.018305131474586
$ws.Range("K20").Value = 1.022689049563897
$ws.Range("L20").Value = 1.017584615332903
$ws.Range("N20").Value = 1.010235939623525
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.011052468389165
$ws.Range("D21").Value = 1.018622785362611
$ws.Range("E21").Value = 1.013386918204838
$ws.Range("I21").Value = 1.025723171786534
$ws.Range("J21").Value = 1.017823115188152
$ws.Range("K21").Value = 1.022268818178707
$ws.Range("L21").Value = 1.017053297395616
$ws.Range("N21").Value = 1.010075700632866
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.010493555301183
$ws.Range("D22").Value = 1.018224991455997
$ws.Range("E22").Value = 1.012918971790797
$ws.Range("I22").Value = 1.025666062032074
$ws.Range("J22").Value = 1.017520058434592
$ws.Range("K22").Value = 1.022004140805096
$ws.Range("L22").Value = 1.016719556440984
$ws.Range("N22").Value = 1.009974926964303
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.010789778119791
$ws.Range("D23").Value = 1.018435840029972
$ws.Range("E23").Value = 1.013166937994047
$ws.Range("I23").Value = 1.025696466669385
$ws.Range("J23").Value = 1.017680713510636
$ws.Range("K23").Value = 1.022144494505991
$ws.Range("L23").Value = 1.016896448017135
$ws.Range("N23").Value = 1.010028351201298
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.011957467662858
$ws.Range("D24").Value = 1.019266573514098
$ws.Range("E24").Value = 1.014145386811302
$ws.Range("I24").Value = 1.025813286960862
$ws.Range("J24").Value = 1.018313212538801
$ws.Range("K24").Value = 1.022696086868942
$ws.Range("L24").Value = 1.017593528254874
$ws.Range("N24").Value = 1.010238625589835
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.013315518509072
$ws.Range("D25").Value = 1.020231847239014
$ws.Range("E25").Value = 1.015285363463984
$ws.Range("I25").Value = 1.025942858082052
$ws.Range("J25").Value = 1.019047171027472
$ws.Range("K25").Value = 1.023334120254233
$ws.Range("L25").Value = 1.018403791573574
$ws.Range("N25").Value = 1.010482511446143
